# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps and the
# Priority value for the two localization-bundle rows (1ebd7b39... and
# b04c7b82...), matching a newly regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" for the two rows
#     that were regenerated (rows 2 and 4) ---
$wsOverview.Range("G2").Value = "2016-09-01 06:18:48"
$wsOverview.Range("G4").Value = "2016-09-01 06:18:48"

# --- zh-cn sheet (rows 2 and 4) ---
# Priority: ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-01 06:18:43"
$wsZhCn.Range("H4").Value = "2016-09-01 06:18:43"
# Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-09-01 06:19:07"
$wsZhCn.Range("K4").Value = "2016-09-01 06:19:07"

# --- de-de sheet (rows 2 and 4) ---
# Priority: ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime (shares the same string as Overview G2/G4)
$wsDeDe.Range("H2").Value = "2016-09-01 06:18:48"
$wsDeDe.Range("H4").Value = "2016-09-01 06:18:48"
# Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-01 06:19:16"
$wsDeDe.Range("K4").Value = "2016-09-01 06:19:16"
